$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($table, $row, $col, $oldVal, $newVal) {
    $cell = $table.Cell($row, $col)
    $current = $cell.Range.Text
    # Cell range text includes trailing cell-mark characters; compare by prefix.
    if (-not $current.StartsWith($oldVal)) {
        throw "Cell($row,$col) expected '$oldVal' but found '$current'"
    }
    $cell.Range.Text = $newVal
}

# Row 1
Set-CellValue $t 1 1 "96÷3=32, 0" "98÷9=10, 8"
Set-CellValue $t 1 2 "23÷8=2, 7" "26÷9=2, 8"
Set-CellValue $t 1 3 "74÷5=14, 4" "76÷6=12, 4"
Set-CellValue $t 1 4 "58÷5=11, 3" "89÷6=14, 5"
Set-CellValue $t 1 5 "95÷5=19, 0" "39÷7=5, 4"

# Row 5
Set-CellValue $t 5 1 "13÷9=1, 4" "61÷8=7, 5"
Set-CellValue $t 5 2 "94÷8=11, 6" "22÷4=5, 2"
Set-CellValue $t 5 3 "70÷8=8, 6" "96÷7=13, 5"
Set-CellValue $t 5 4 "37÷7=5, 2" "34÷9=3, 7"
Set-CellValue $t 5 5 "95÷3=31, 2" "93÷4=23, 1"

# Row 9
Set-CellValue $t 9 1 "43÷9=4, 7" "72÷2=36, 0"
Set-CellValue $t 9 2 "24÷9=2, 6" "27÷9=3, 0"
Set-CellValue $t 9 3 "19÷2=9, 1" "96÷3=32, 0"
Set-CellValue $t 9 4 "95÷2=47, 1" "55÷5=11, 0"
Set-CellValue $t 9 5 "58÷8=7, 2" "90÷7=12, 6"

# Row 13
Set-CellValue $t 13 1 "85÷5=17, 0" "17÷3=5, 2"
Set-CellValue $t 13 2 "87÷8=10, 7" "52÷8=6, 4"
Set-CellValue $t 13 3 "17÷8=2, 1" "28÷5=5, 3"
Set-CellValue $t 13 4 "79÷2=39, 1" "13÷9=1, 4"
Set-CellValue $t 13 5 "84÷7=12, 0" "68÷6=11, 2"

# Row 17
Set-CellValue $t 17 1 "43÷9=4, 7" "92÷8=11, 4"
Set-CellValue $t 17 2 "54÷6=9, 0" "89÷2=44, 1"
Set-CellValue $t 17 3 "58÷3=19, 1" "84÷3=28, 0"
Set-CellValue $t 17 4 "25÷3=8, 1" "52÷5=10, 2"
Set-CellValue $t 17 5 "29÷9=3, 2" "11÷5=2, 1"

Write-Output "Updated 25 table cells successfully"
